# ------------------------------------------------------------------
# "zone statistics and cfg statistics refactored for case if all
# devices are local" -- adds a new 'CKR' project entry to the
# 'report' sheet (a new column inserted before the existing data),
# flips a handful of zoning-statistics flag cells on 'service_tables'
# and appends one new service_tables row describing a new
# 'zoning_modified' step, and widens the hidden _FilterDatabase
# defined name to include that new row.
# ------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # report
$ws2 = $wb.Worksheets.Item(2)   # service_tables

# ===================================================================
# 1) "report" sheet: insert a new column for the "CKR" project
# ===================================================================

# Inserting at column C shifts the existing C:AA project columns one
# slot to the right (D:AB), matching every other column's width/style
# along with it.
$ws1.Columns.Item(3).Insert()

# Row 1 is a "name"/"value" mini-header that lived at A1/E1; after the
# shift the "value" label sits at F1, but it belongs at C1 (directly
# above the new project's data column).
$valueLabel = $ws1.Range("F1").Value()
$ws1.Range("C1").Value = $valueLabel
$ws1.Range("A1").Copy()
$ws1.Range("C1").PasteSpecial(-4122)
$ws1.Range("F1").Clear()

# Populate the new column with the CKR project's details.
$ws1.Range("C2").Value = "CKR"

$ws1.Range("C3").Value = 44244
$ws1.Range("D3").Copy()
$ws1.Range("C3").PasteSpecial(-4122)

$ws1.Range("C4").Value = "SAN Implementation"

$ws1.Range("C5").Value = "C:\Users\vlasenko\Documents\01.CUSTOMERS\CKR\SAN Impl FEB2021"
$ws1.Range("C6").Value = "C:\Users\vlasenko\Documents\06.CONFIGS\CKR\FEB21"

$ws1.Range("E14").Select()

# ===================================================================
# 2) "service_tables" sheet: zone/cfg statistics rows now report
#    "local-only" status directly (G flips) instead of depending on a
#    separate local/overall split (H56/G71 go back to 0)
# ===================================================================

$ws2.Range("H56").Value = 0
$ws2.Range("G71").Value = 0

$localOnlyRows = 82,83,84,85,86,87,88,89,91,93,96
foreach ($r in $localOnlyRows) {
    $ws2.Range("G$r").Value = 1
}

# ------------------------------------------------------------------
# New row describing the "zoning_modified" analysis step, appended
# right after the existing last row (96).
# ------------------------------------------------------------------
$ws2.Range("B97").Value = "analysis_zoning"

$ws2.Range("D97").Value = "analysis"
$ws2.Range("D92").Copy()
$ws2.Range("D97").PasteSpecial(-4122)

$ws2.Range("F97").Value = "zoning_modified"
$ws2.Range("F92").Copy()
$ws2.Range("F97").PasteSpecial(-4122)

$ws2.Range("G97").Value = 0
$ws2.Range("H97").Value = 0
$ws2.Range("I97").Value = "DATA ANALYSIS 6. ZONING CONFIGURATION"

$ws2.Range("J97").Value = "-"
$ws2.Range("J92").Copy()
$ws2.Range("J97").PasteSpecial(-4122)

$ws2.Activate()
$ws2.Range("G80").Select()

# ===================================================================
# 3) Widen the hidden _FilterDatabase defined name to cover the new
#    row 97 on service_tables.
# ===================================================================
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $nm = $wb.Names.Item($i)
    if ($nm.Name -like "*_FilterDatabase*") {
        $nm.RefersTo = "=service_tables!`$A`$1:`$J`$97"
    }
}
